$wb = $excel.ActiveWorkbook

# --- Update the Tracker sheet view: drop the "tabSelected" focus (another sheet
# becomes active below), scroll the frozen pane back to the top, and select A1:G1 ---
$tracker = $wb.Worksheets.Item("Tracker")
[void]$tracker.Range("A1:G1").Select()

# --- Insert the new "Antech-Issues" sheet right after "MockTest" (and before
# "WorkDays"), matching the sheet order / sheetId seen in the target workbook ---
$afterSheet = $wb.Worksheets.Item("MockTest")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "Antech-Issues"

# Header row (reuses the same column headers / bold styling as the Tracker sheet)
$newSheet.Range("A1").Value = "Item Number"
$newSheet.Range("B1").Value = "Item Description"
$newSheet.Range("C1").Value = "Type"
$newSheet.Range("D1").Value = "Owned by"
$newSheet.Range("E1").Value = "Priority"
$newSheet.Range("F1").Value = "Status"
$newSheet.Range("G1").Value = "Comments"
$newSheet.Range("A1:G1").Font.Bold = $true
$newSheet.Range("A1").HorizontalAlignment = -4108

# Antech issue rows
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = "The outside photo (360 degree) of our institute should be proper with banners and no blury images"
$newSheet.Range("C2").Value = "Marketing"
$newSheet.Range("D2").Value = "Antech"
$newSheet.Range("E2").Value = 5
$newSheet.Range("F2").Value = "Todo"

$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = "Address of our institute should be : N 1/25, Kunal Road,Patuli, Kolkata - 700094"
$newSheet.Range("C3").Value = "Marketing"
$newSheet.Range("D3").Value = "Antech"
$newSheet.Range("E3").Value = 5
$newSheet.Range("F3").Value = "Todo"

$newSheet.Range("A4").Value = 3
$newSheet.Range("B4").Value = "Opening Hours should be 10 AM to 7 PM."
$newSheet.Range("C4").Value = "Marketing"
$newSheet.Range("D4").Value = "Antech"
$newSheet.Range("E4").Value = 5
$newSheet.Range("F4").Value = "Todo"

$newSheet.Range("A5").Value = 4
$newSheet.Range("B5").Value = "In Appointments: wa.me area, there should be whats app icon."
$newSheet.Range("C5").Value = "Marketing"
$newSheet.Range("D5").Value = "Antech"
$newSheet.Range("E5").Value = 5
$newSheet.Range("F5").Value = "Todo"

$newSheet.Range("A6").Value = 5
$newSheet.Range("B6").Value = "There is no place for giving review comments"
$newSheet.Range("C6").Value = "Marketing"
$newSheet.Range("D6").Value = "Antech"
$newSheet.Range("E6").Value = 5
$newSheet.Range("F6").Value = "Todo"

$newSheet.Range("A7").Value = 6
$newSheet.Range("B7").Value = "On pressing Website button the page should redirect to google.anodiam.mybusiness.site"
$newSheet.Range("C7").Value = "Marketing"
$newSheet.Range("D7").Value = "Antech"
$newSheet.Range("E7").Value = 5
$newSheet.Range("F7").Value = "Todo"

# Column widths, roughly matching the authored sheet
$newSheet.Columns.Item(1).ColumnWidth = 12.22
$newSheet.Columns.Item(2).ColumnWidth = 82.33
$newSheet.Columns.Item(3).ColumnWidth = 19.22
$newSheet.Columns.Item(4).ColumnWidth = 12.33
$newSheet.Columns.Item(7).ColumnWidth = 10.11

# Item Number column is centered, same as on the Tracker sheet
$newSheet.Columns.Item(1).HorizontalAlignment = -4108

# Data validation drop-downs for Type / Owned by / Status columns
$newSheet.Range("F1:F1048576").Validation.Add(3, 1, 3, '"Todo,WIP,Done"')
$newSheet.Range("C1:C1048576").Validation.Add(3, 1, 3, '"Marketing"')
$newSheet.Range("D1:D1048576").Validation.Add(3, 1, 3, '"Antech,Anirban,Debashish,Sayan,Rahul"')

# The new sheet becomes the active tab (matches activeTab=2 / tabSelected on the
# new sheet in the target workbook)
[void]$newSheet.Range("B2").Select()
[void]$newSheet.Activate()
